$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected FilesTab Neo4j Cypher script (B4): drop the `File Type` and
# `Breed` columns from the RETURN clause.
$newFilesQuery = "MATCH (f:file)-->(parent)`n" +
  "WITH DISTINCT f, parent`n" +
  "MATCH (f)-[*]->(c:case)<--(demo:demographic)`n" +
  "WHERE demo.breed IN ['Weimaraner']`n" +
  "OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n" +
  "OPTIONAL MATCH (samp:sample)-->(c)`n" +
  "WITH DISTINCT f, parent, c, demo, diag, s`n" +
  "RETURN  coalesce(f.file_name, '') AS ``File Name``,`n" +
  "        coalesce(labels(parent)[0], '') AS ``Association``,`n" +
  "        coalesce(f.file_description, '') AS ``Description``,`n" +
  "        coalesce(f.file_format, '') AS ``Format``,`n" +
  "        coalesce(f.file_size, '') AS ``Size``,`n" +
  "        coalesce(c.case_id, '') AS ``Case ID``,`n" +
  "        coalesce(diag.disease_term,'') AS Diagnosis , `n" +
  "        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFilesQuery

# The shorter script text wraps to fewer lines, so the row shrinks.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection/scroll moved onto the corrected Files script cell.
$ws.Activate()
$ws.Range("B4").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 4
$window.ScrollColumn = 1
